$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), reusing the existing header
# style (copy format+value from H1, then overwrite with the new text).
$ws.Range("H1").Copy($ws.Range("I1:J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Values for columns I (I0) and J (IF) for rows 2 through 17
$colI = @(8, 8, 8, 7, 3, 6, 6, 6, 6, 3, 4, 4, 4, 6, 5, 6)
$colJ = @(9, 9, 8, 8, 6, 7, 6, 6, 7, 3, 5, 4, 4, 6, 6, 6)

for ($i = 0; $i -lt $colI.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $colI[$i]
    $ws.Cells.Item($row, 10).Value = $colJ[$i]
}
